# ------------------------------------------------------------------
# Add a new "2022-Q1" sheet with fund-holding detail data (placed
# right before the "总计" summary sheet) and prepend a matching
# summary row to the "总计" sheet.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$templateSheet = $wb.Worksheets.Item("2021-Q4")

$fundRows = @(
    @("000242", "景顺长城策略精选", "20.09", "91.24", "4.98", "1.0005", 4),
    @("166023", "中欧瑞丰灵活配置混合（LOF）A", "32.40", "85.04", "3.04", "0.9850", 10),
    @("010147", "博道嘉兴一年持有期混合", "14.00", "92.27", "5.52", "0.7728", 5),
    @("010967", "博道嘉丰混合A", "19.02", "90.02", "3.84", "0.7304", 6),
    @("008208", "博道嘉泰回报混合", "16.72", "77.97", "3.70", "0.6186", 5),
    @("008467", "博道嘉瑞混合A", "16.14", "77.44", "3.49", "0.5633", 6),
    @("010201", "农银汇理智增一年定期开放混合", "9.82", "62.05", "3.87", "0.3800", 6),
    @("008793", "博道嘉元混合A", "9.21", "76.45", "3.59", "0.3306", 5),
    @("501019", "国泰国证航天军工指数（LOF）", "9.84", "93.43", "3.13", "0.3080", 8),
    @("000314", "招商瑞丰灵活配置混合A", "14.01", "40.96", "1.90", "0.2662", 6),
    @("001427", "招商丰泽灵活配置混合A", "13.94", "41.14", "1.78", "0.2481", 6),
    @("004278", "东方红智逸沪港深定期开放混合", "34.95", "22.66", "0.56", "0.1957", 10),
    @("161123", "易方达并购重组指数（LOF）", "4.78", "94.71", "3.93", "0.1879", 7),
    @("004475", "华泰柏瑞富利灵活配置混合", "4.43", "93.22", "3.49", "0.1546", 7),
    @("519756", "交银施罗德国企改革灵活配置混合", "5.69", "87.70", "2.56", "0.1457", 9),
    @("002389", "招商安德灵活配置混合A", "8.03", "44.29", "1.74", "0.1397", 8),
    @("010968", "博道嘉丰混合C", "3.58", "90.02", "3.84", "0.1375", 6),
    @("004142", "招商盛合灵活配置混合A", "8.42", "43.54", "1.62", "0.1364", 10),
    @("004143", "招商盛合灵活配置混合C", "7.31", "43.54", "1.62", "0.1184", 10),
    @("001656", "农银汇理中国优势灵活配置混合", "2.92", "68.89", "3.91", "0.1142", 6),
    @("002776", "招商安荣灵活配置混合A", "6.50", "41.03", "1.75", "0.1138", 6),
    @("002819", "招商丰美灵活配置混合A", "7.83", "29.10", "1.36", "0.1065", 8),
    @("005412", "金信民长灵活配置混合A", "2.37", "70.53", "4.37", "0.1036", 8),
    @("002581", "招商丰凯灵活配置混合A", "5.76", "40.60", "1.79", "0.1031", 1),
    @("002017", "招商瑞丰灵活配置混合C", "4.66", "40.96", "1.90", "0.0885", 6),
    @("001247", "华泰柏瑞新利灵活配置混合A", "10.68", "20.70", "0.70", "0.0748", 7),
    @("004010", "华泰柏瑞鼎利灵活配置混合A", "10.24", "21.26", "0.70", "0.0717", 7),
    @("005413", "金信民长灵活配置混合C", "1.30", "70.53", "4.37", "0.0568", 8),
    @("008794", "博道嘉元混合C", "1.58", "76.45", "3.59", "0.0567", 5),
    @("002390", "招商安德灵活配置混合C", "2.55", "44.29", "1.74", "0.0444", 8),
    @("002091", "华泰柏瑞新利灵活配置混合C", "5.98", "20.70", "0.70", "0.0419", 7),
    @("008468", "博道嘉瑞混合C", "1.13", "77.44", "3.49", "0.0394", 6),
    @("004740", "中欧瑞丰灵活配置混合（LOF）C", "1.28", "85.04", "3.04", "0.0389", 10),
    @("001446", "招商丰泽灵活配置混合C", "1.99", "41.14", "1.78", "0.0354", 6),
    @("002582", "招商丰凯灵活配置混合C", "1.46", "40.60", "1.79", "0.0261", 1),
    @("006377", "广发趋势动力灵活配置混合", "0.82", "90.19", "3.05", "0.0250", 9),
    @("004011", "华泰柏瑞鼎利灵活配置混合C", "3.47", "21.26", "0.70", "0.0243", 7),
    @("002820", "招商丰美灵活配置混合C", "0.73", "29.10", "1.36", "0.0099", 8),
    @("002777", "招商安荣灵活配置混合C", "0.54", "41.03", "1.75", "0.0094", 6),
    @("410006", "华富策略精选混合", "0.11", "77.80", "2.75", "0.0030", 7),
    @("006535", "恒生前海恒锦裕利混合A", "0.14", "29.05", "1.26", "0.0018", 6),
    @("006536", "恒生前海恒锦裕利混合C", "0.11", "29.05", "1.26", "0.0014", 6)
)

# --- Create the new "2022-Q1" worksheet, positioned right before "总计" ---
# (NOTE: the sheet reference passed as "Before" gets repointed to the
#  freshly-created sheet once it is renamed, so the "总计" sheet must be
#  looked up again afterwards instead of reusing that variable.)
$beforeRef = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($beforeRef)
$newSheet.Name = "2022-Q1"

$totalSheet = $wb.Worksheets.Item("总计")

# --- Header row (row 1): B..H, matching the layout used by the other
#     quarterly sheets (基金代码 / 基金名称 / 基金规模 / 股票总仓位 / 仓位占比 / 持有市值(亿元) / 仓位排名)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# copy header-cell styling (bold, centered, bordered) from an existing sheet
$templateSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# --- Data rows (2..43) ---
for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# copy the index-column styling (bold, centered, bordered) used by column A
# of the other quarterly sheets onto the new sheet's A2:A43 range
$templateSheet.Range("A2").Copy()
$newSheet.Range("A2:A43").PasteSpecial(-4122)

# --- Update the "总计" (summary) sheet: insert a new row for 2022-Q1 on top ---
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 42
$totalSheet.Range("D2").Value = 8.609999999999999

# restore the index-column style (s=2) on the new A2 cell to match A3:A5
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
